$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D13: <FD> void-returning function declaration semantic rule
# "push <SB>.innerVarAmount vars from var symbol table" line removed
$ws.Range("D13").Value = "<FD>.returnType != <SB>.rT: ERROR`n<FD>.IR = 'def func(' + <FP>.IR + ')' + <SB>.IR + '}'`n<FD>.returnType = 'void'`n<FD>.paramName = <FP>.paramName`n<FD>.paramType = <FP>.paramType`ninsert into func symbol table"

# D11: <FD> int-returning function declaration semantic rule
$ws.Range("D11").Value = "<FD>.returnType != <SB>.rT: ERROR`n<FD>.IR = 'def func(' + <FP>.IR + ')' + <SB>.IR + '}'`n<FD>.returnType = 'int'`n<FD>.paramName = <FP>.paramName`n<FD>.paramType = <FP>.paramType`ninsert into func symbol table"

# D12: <FD> float-returning function declaration semantic rule
$ws.Range("D12").Value = "<FD>.returnType != <SB>.rT: ERROR`n<FD>.IR = 'def func(' + <FP>.IR + ')' + <SB>.IR + '}'`n<FD>.returnType = 'float'`n<FD>.paramName = <FP>.paramName`n<FD>.paramType = <FP>.paramType`ninsert into func symbol table"

# D21: <SB> statement block semantic rule
# dropped "<SB>.innerVarAmount = <S>.vA" line and changed pop line to reference <S>.innerVarAmount
$ws.Range("D21").Value = "<SB>.IR = <S>.IR`n<SB>.returnType = <S>.returnType`npop <S>.innerVarAmount vars from var symbol table"

# D22: <Stmts1> statement-list semantic rule
# fixed typos: <Stmt>vA -> <Stmt>.iVA and <Stmts2>.vA -> <Stmts2>.iVA
$ws.Range("D22").Value = "<Stmts1>.IR = <Stmt>.IR + <Stmts2>.IR`nif <S> <Stmts2> rT equal: <Stmts1>.rT = <Stmt>.rT`nelse: ERROR`n<Stmts1>.innerVarAmount = <Stmt>.iVA + <Stmts2>.iVA"

# D31: <RS> void return-statement semantic rule
# changed generated IR text from 'ret ;' to 'ret void'
$ws.Range("D31").Value = "<RS>.IR = 'ret void'`n<RS>.returnType = 'void'"
